# 23 Nov 2021 3rd commit
# Turn the CCPA URL "status" flag from ON to OFF for every URL row except
# the first block (rows 2-9), i.e. set C10:C92 on the "CCPAUrls" sheet to
# "OFF" (this introduces a brand-new shared string "OFF"), and move the
# selection to the cell range that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCPAUrls")
$ws.Activate()

# Flip the status column for rows 10 through 92 to "OFF".
$ws.Range("C10:C92").Value = "OFF"

# Match the author's final selection/view: C10 active, C10:C92 selected,
# scrolled down so row 72 is the first visible row under the frozen header.
$ws.Range("C10:C92").Select()
$excel.ActiveWindow.ScrollRow = 72
